$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 427; this shifts rows 427..457 down to 428..458
$ws.Rows.Item(427).Insert()

# Populate the new row 427 with the new weekly record
$ws.Cells.Item(427, 1).Value = 4
$ws.Cells.Item(427, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(427, 3).Value = "Los Lagos"
$ws.Cells.Item(427, 4).Value = 45013
$ws.Cells.Item(427, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(427, 5).Value = 10
$ws.Cells.Item(427, 6).Value = 100114014
$ws.Cells.Item(427, 7).Value = "Betarraga"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 1200
$ws.Cells.Item(427, 11).Value = 1100
$ws.Cells.Item(427, 12).Value = 1100
$ws.Cells.Item(427, 13).Value = 1100
$ws.Cells.Item(427, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(427, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(427, 16).Value = 220
$ws.Cells.Item(427, 17).Value = 5
$ws.Cells.Item(427, 18).Value = "Hortaliza"
